$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row (row 33) with data for "Machine 32"
$ws.Cells.Item(33, 1).Value = 10032
$ws.Cells.Item(33, 2).Value = "Machine 32"
$ws.Cells.Item(33, 3).Value = "F4-30-B9-D4-CD-6F"
$ws.Cells.Item(33, 4).Value = "FB5962911665"
$ws.Cells.Item(33, 5).Value = "192.168.0.358"
$ws.Cells.Item(33, 6).Value = 1001
$ws.Cells.Item(33, 7).Value = "eng"
$ws.Cells.Item(33, 8).Value = $true
$ws.Cells.Item(33, 9).Value = "superadmin"
$ws.Cells.Item(33, 10).Value = "now()"

# Update the active cell selection to match target state
$ws.Range("J29").Select()
